$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet (Version, Date, Contact) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Duplicate "Include from FSIII" into a new "Include from FSIII 2" sheet ---
$src = $wb.Worksheets.Item("Include from FSIII")
$src.Copy($null, $src)
$copy = $wb.Worksheets.Item($wb.Worksheets.Count)
$copy.Name = "Include from FSIII 2"

# --- Update the original sheet's code value to the new UUID, keep the copy's old "A" value ---
$src.Range("C2").Value = "95851822-5a33-4349-a1f2-9b1245369bf5"
